$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.872.50'
$ws.Range('E2').Value = '  +0.23%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.664.57'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.64'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.73'
$ws.Range('E6').Value = '  +1.18%  '
$ws.Range('E7').Value = '  +4.90%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -2.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.87'
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('E12').Value = '  +1.69%  '
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('E14').Value = '  -1.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.143.39'
$ws.Range('E15').Value = '  -0.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.718.59'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.657.65'
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.62'
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.81'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -1.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '351.93'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.92'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.82'
$ws.Range('E24').Value = '  +11.03%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '572.43'
$ws.Range('E28').Value = '  +8.20%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.19'
$ws.Range('E29').Value = '  +1.26%  '
$ws.Range('E30').Value = '  -2.72%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('E33').Value = '  +2.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.76'
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.56'
$ws.Range('E35').Value = '  +1.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.424'
$ws.Range('E36').Value = '  -0.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.62'
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.97'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '154.50'
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '161.69'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.12'
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('B44').Value = 'Hedera'
$ws.Range('C44').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0620'
$ws.Range('E44').Value = '  +1.33%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.32'
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.22'
$ws.Range('E46').Value = '  +1.28%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.645'
$ws.Range('E47').Value = '  +0.36%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0258'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.103'
$ws.Range('E49').Value = '  +2.00%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.81'
$ws.Range('E50').Value = '  -1.85%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0246'
$ws.Range('E51').Value = '  -8.40%  '
